# Add row 19 to sheet1: new SmartScore submission from Andrew Moody
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A19").Value = "Andrew Moody_20251202_125927"
$ws.Range("C19").Value = "Andrew Moody"
$ws.Range("D19").Value = 18
$ws.Range("E19").Value = "Male"
$ws.Range("F19").Value = "2025-12-02 12:59:27"
$ws.Range("G19").Value = @"
{
  "portion": 0.6,
  "diet": 0.7142857142857143,
  "salt": 0.2,
  "fat": 0.8,
  "natural": 0.2,
  "convenience": 0.8,
  "price": 1.0
}
"@
$ws.Range("H19").Value = "Maruchan Ramen Sabor Pollo"
$ws.Range("I19").Value = "'0.572"
$ws.Range("J19").Value = "Sabor clásico, económico, alto en sodio, no saludable, nostálgico"
$ws.Range("K19").Value = "Nongshim Neoguri Spicy Seafood"
$ws.Range("L19").Value = "'0.554"
$ws.Range("M19").Value = "Sabor a marisco, umami, picante equilibrado, buena textura, algo salado"
$ws.Range("N19").Value = "Nissin Chow Mein Teriyaki Beef"
$ws.Range("O19").Value = "'0.540"
$ws.Range("P19").Value = "Fácil de preparar, porción generosa, salsa suave, necesita mejoras, alto en grasa"
$ws.Range("Q19").Value = "Velveeta Original Shells & Cheese (microwave cups)"
$ws.Range("R19").Value = "'0.591"
$ws.Range("S19").Value = "Muy cremoso, porción individual, rápido, salado, ideal para niños"
$ws.Range("T19").Value = "Kraft Macaroni & Cheese Dinner"
$ws.Range("U19").Value = "'0.541"
$ws.Range("V19").Value = "Sabor nostálgico, clásico americano, fácil, no muy nutritivo, barato"
$ws.Range("W19").Value = "Annie’s Shells & White Cheddar"
$ws.Range("X19").Value = "'0.460"
$ws.Range("Y19").Value = "Queso blanco real, sin colorantes, sabor casero, menos salado, buena para niños"
$ws.Range("Z19").Value = "StarKist Chicken Creations (Chicken Salad)"
$ws.Range("AA19").Value = "'0.674"
$ws.Range("AB19").Value = "Portátil, saludable, fácil, buena textura, sabor suave"
$ws.Range("AC19").Value = "Jack Link’s Beef Jerky Original"
$ws.Range("AD19").Value = "'0.659"
$ws.Range("AE19").Value = "Ahumado, sabroso, alto en proteína, snack ideal, porción pequeña"
$ws.Range("AF19").Value = "Wild Planet Wild Tuna Pasta Salad"
$ws.Range("AG19").Value = "'0.650"
$ws.Range("AH19").Value = "Sabor fresco, buena proteína, saludable, porción algo pequeña"

# Reset row height back to automatic so no explicit customHeight is stored
$ws.Rows(19).AutoFit() | Out-Null

